$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after the header row; Excel shifts the existing
# data (old rows 2-19) down to rows 4-21, carrying their old A/B values too.
$ws.Rows("2:3").Insert()

# The insert copies formatting down from the header row (bold/border); the
# new rows should be plain data rows like the others, so clear that format.
$ws.Rows("2:3").ClearFormats()

# Populate the new accelerometer/gyroscope readings for the two inserted rows
$ws.Cells.Item(2,3).Value2 = -1.97176456451416
$ws.Cells.Item(2,4).Value2 = 1.745009422302246
$ws.Cells.Item(2,5).Value2 = 0.4838592410087585
$ws.Cells.Item(2,6).Value2 = 0.1914996167887811
$ws.Cells.Item(2,7).Value2 = 0.03030422819859344
$ws.Cells.Item(2,8).Value2 = 0.02057685541069637

$ws.Cells.Item(3,3).Value2 = -1.845728397369385
$ws.Cells.Item(3,4).Value2 = 1.672563552856445
$ws.Cells.Item(3,5).Value2 = 0.5211508870124817
$ws.Cells.Item(3,6).Value2 = 0.1882859338884768
$ws.Cells.Item(3,7).Value2 = 0.0655750582480559
$ws.Cells.Item(3,8).Value2 = 0.4253946024438608

# timestamp (A) and label (B) follow a fixed, regular sequence across all 20
# data rows (0, 100, 200, ... ms; always "falling") -- reset them for every
# row so they read correctly after the shift above.
For ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r,1).Value2 = ($r - 2) * 100
    $ws.Cells.Item($r,2).Value2 = "falling"
}
